# Revision: Add Tutor Management and Exam History Features
# Replace the "Nisn" (NISN) identifier column with a "Nik" (NIK) column
# holding longer text identifiers, widen column A to fit, and move the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Nisn -> Nik
$ws.Range("A1").Value = "Nik"

# NIK values replace the old numeric NISN values (101..105).
# Store as text (left-aligned "49" text format, same as the Phone column)
# so the long digit strings are not turned into floating point numbers.
$ws.Range("A2:A6").NumberFormat = "@"

$nikValues = @("3518041000000001", "3518041000000002", "3518041000000003", "3518041000000004", "3518041000000005")
for ($i = 0; $i -lt $nikValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $nikValues[$i]
}

# Widen column A to fit the new NIK values (~26.57 characters).
$ws.Columns.Item(1).ColumnWidth = 25.65

# Move the active selection.
$ws.Range("I12").Select()
